# 224614 Add ERP ID update for product sync command
#
# Replace the old ERP/product sync id "1213-3316" with the new id
# "6199-6365" across every sheet that references it, update the
# "Action" column on the Items sheet to "publish", drop the stray
# "Error" column (S) on the Items sheet, and refresh the active
# selections that the workbook was left on.

$wb = $excel.ActiveWorkbook

# --- General -------------------------------------------------------------
$ws = $wb.Worksheets.Item("General")
$ws.Range("B3").Value = "PRD-6199-6365"

# --- Parameters Groups -----------------------------------------------------
$ws = $wb.Worksheets.Item("Parameters Groups")
$ws.Range("A2").Value = "PGR-6199-6365-0002"
$ws.Range("A3").Value = "PGR-6199-6365-0003"

# --- Items Groups ------------------------------------------------------
$ws = $wb.Worksheets.Item("Items Groups")
$ws.Range("A2").Value = "IGR-6199-6365-0002"
$ws.Range("A3").Value = "IGR-6199-6365-0003"
$ws.Activate()
$ws.Range("A2").Select()

# --- Agreements Parameters -----------------------------------------------
$ws = $wb.Worksheets.Item("Agreements Parameters")
$ws.Range("A2").Value = "PAR-6199-6365-0001"
$ws.Range("H2").Value = "PGR-6199-6365-0002"
$ws.Range("A3").Value = "PAR-6199-6365-0002"

# --- Item Parameters -------------------------------------------------------
$ws = $wb.Worksheets.Item("Item Parameters")
$ws.Range("A2").Value = "PAR-6199-6365-0003"
$ws.Range("A3").Value = "PAR-6199-6365-0004"

# --- Request Parameters ---------------------------------------------------
$ws = $wb.Worksheets.Item("Request Parameters")
$ws.Range("A2").Value = "PAR-6199-6365-0005"
$ws.Range("A3").Value = "PAR-6199-6365-0006"

# --- Subscription Parameters -----------------------------------------------
$ws = $wb.Worksheets.Item("Subscription Parameters")
$ws.Range("A2").Value = "PAR-6199-6365-0007"
$ws.Range("A3").Value = "PAR-6199-6365-0008"

# --- Items -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Items")
$ws.Range("S1").ClearContents()
$ws.Range("A2").Value = "ITM-6199-6365-0001"
$ws.Range("C2").Value = "publish"
$ws.Range("J2").Value = "IGR-6199-6365-0002"
$ws.Range("A3").Value = "ITM-6199-6365-0002"
$ws.Range("C3").Value = "publish"
$ws.Range("J3").Value = "IGR-6199-6365-0002"
$ws.Activate()
$ws.Range("C16").Select()

# --- Templates ---------------------------------------------------------
$ws = $wb.Worksheets.Item("Templates")
$ws.Range("A2").Value = "TPL-6199-6365-0005"
$ws.Range("F2").Value = "Test content **Azure** {{ PAR-6199-6365-0001 }}"
$ws.Range("A3").Value = "TPL-6199-6365-0006"
